$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates for the crypto price/volume refresh
$ws.Range("D2").Value = '23.492.28'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.650.37'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9997'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '300.30'
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3787'
$ws.Range("E7").Value = '  -1.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.65'
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3508'
$ws.Range("E9").Value = '  -2.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.224'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08056'
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.10'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.318'
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.267'
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("D17").Value = '1.650.08'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.32'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06982'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.634'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.46'
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.46'
$ws.Range("E23").Value = '  -1.68%  '
$ws.Range("D24").Value = '23.497.50'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("E25").Value = '  -3.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.025'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.07'
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.72'
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.176'
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.94'
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").Value = '1.839.89'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.889'
$ws.Range("E32").Value = '  -4.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.136'
$ws.Range("E33").Value = '  -5.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.21'
$ws.Range("E34").Value = '  -7.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9927'
$ws.Range("E35").Value = '  -6.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02716'
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08807'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.941'
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2427'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06827'
$ws.Range("E40").Value = '  -2.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.91'
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6902'
$ws.Range("E42").Value = '  -1.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.299'
$ws.Range("E43").Value = '  -2.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.61'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9988'
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6399'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.248'
$ws.Range("E47").Value = '  -2.47%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07693'
$ws.Range("E49").Value = '  -2.65%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.14'
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.235'
$ws.Range("E51").Value = '  +2.73%  '
